$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.164.55"
$ws.Range("E2").Value = "  -3.57%  "

$ws.Range("D3").Value = "2.973.23"
$ws.Range("E3").Value = "  -0.53%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "557.43"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.75%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.96"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.61%  "

$ws.Range("E7").Value = "  -0.09%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.519"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.58%  "

$ws.Range("D9").Value = "2.967.73"
$ws.Range("E9").Value = "  -0.58%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.129"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.56%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "4.88"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.23%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.452"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.67%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000224"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.62%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.13"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.80%  "

$ws.Range("E15").Value = "  +0.77%  "

$ws.Range("D16").Value = "3.453.31"
$ws.Range("E16").Value = "  -0.58%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.82"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +9.71%  "

$ws.Range("D18").Value = "2.951.94"
$ws.Range("E18").Value = "  -1.04%  "

$ws.Range("D19").Value = "58.016.37"
$ws.Range("E19").Value = "  -3.57%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "420.13"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.80%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.23"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.91%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.688"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.92%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.01"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.20%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.10"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.36%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "79.69"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.75%  "

$ws.Range("E26").Value = "  +0.02%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.09%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.50"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.90%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.61"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.48%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.90%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.35"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.08%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.10"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.71%  "

$ws.Range("E33").Value = "  +8.29%  "

$ws.Range("B34").Value = "Stacks"
$ws.Range("C34").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.14"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.48%  "

$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.68"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.55%  "

$ws.Range("B36").Value = "Mantle"
$ws.Range("C36").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.944"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.99%  "

$ws.Range("B37").Value = "PEPE"
$ws.Range("C37").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D37").Value = "0.0₃0701"
$ws.Range("E37").Value = "  +6.44%  "

$ws.Range("B38").Value = "OKB"
$ws.Range("C38").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "48.71"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.83%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.49"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.95%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.59"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.72%  "

$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0352"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.82%  "

$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.109"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.05%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "380.83"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.32%  "

$ws.Range("D44").Value = "2.677.49"
$ws.Range("E44").Value = "  +1.67%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.243"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.29%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "122.82"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.93%  "

$ws.Range("E48").Value = "  +2.82%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.01"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.06%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.63"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.54%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.03"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.56%  "
